$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new Location column (E) - written in the same order the source
# workbook's shared-string table was built in
$ws.Range("E3").Value = "SouthCentralUS"
$ws.Range("E1").Value = "Location"
$ws.Range("E2").Value = "EastUS"

# Update existing ResourceGroupName values (column B)
$ws.Range("B2").Value = "rg-hpc-azhop-deploy"
$ws.Range("B3").Value = "rg-dev-mg1311-jumpbox"

# Update existing TagValue values (column D)
$ws.Range("D2").Value = "Value07"
$ws.Range("D3").Value = "Value07"

# Match the new column E width / best-fit formatting (closest attainable
# value to the source file's computed best-fit width of 13.88671875)
$ws.Columns.Item(5).ColumnWidth = 13

# Move the active selection to D4, matching the post-edit state
$ws.Range("D4").Select()
